$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column F ("dSF") per row, per the repull/recalculation.
$updates = @{
    3  = -2
    4  = -4
    5  = -1
    6  = 2
    7  = 2
    8  = 1
    9  = -2
    10 = -2
    12 = 3
    13 = -2
    14 = -1
    15 = 5
    16 = 1
    17 = -2
    18 = -4
    19 = 1
    20 = -2
    22 = 2
    23 = -2
    24 = -2
    25 = -4
    26 = 6
    27 = -4
    28 = -6
    29 = -6
    31 = 2
    32 = 1
    33 = 3
    34 = -1
    35 = 10
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
